# Auto-generated edit script: updates cryptos list per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '63.844.48'
$ws.Range('E2').Value = '  -5.19%  '
$ws.Range('D3').Value = '3.282.11'
$ws.Range('E3').Value = '  -6.38%  '
$c = $ws.Range('D4')
$c.Value = '''1.00'
$c.Style = "Normal"
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '175.34'
$ws.Range('E5').Value = '  -12.42%  '
$c = $ws.Range('D6')
$c.Value = '''522.39'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -5.45%  '
$c = $ws.Range('D7')
$c.Value = '''0.603'
$c.Style = "Normal"
$ws.Range('E7').Value = '  -0.74%  '
$ws.Range('D8').Value = '3.272.93'
$ws.Range('E8').Value = '  -6.37%  '
$ws.Range('E9').Value = '  +0.05%  '
$c = $ws.Range('D10')
$c.Value = '''0.603'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -7.79%  '
$ws.Range('D11').Value = '56.82'
$ws.Range('E11').Value = '  -8.63%  '
$c = $ws.Range('D12')
$c.Value = '''0.132'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -7.60%  '
$c = $ws.Range('D13')
$c.Value = '''0.0000257'
$c.Style = "Normal"
$ws.Range('E13').Value = '  -4.93%  '
$c = $ws.Range('D14')
$c.Value = '''9.02'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -8.09%  '
$ws.Range('D15').Value = '3.814.24'
$ws.Range('E15').Value = '  -6.24%  '
$ws.Range('D16').Value = '3.284.88'
$ws.Range('E16').Value = '  -6.33%  '
$c = $ws.Range('D17')
$c.Value = '''0.116'
$c.Style = "Normal"
$ws.Range('E17').Value = '  -5.98%  '
$ws.Range('D18').Value = '63.827.05'
$ws.Range('E18').Value = '  -4.87%  '
$c = $ws.Range('D19')
$c.Value = '''17.33'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -5.99%  '
$c = $ws.Range('D20')
$c.Value = '''11.00'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -6.69%  '
$c = $ws.Range('D21')
$c.Value = '''0.949'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -7.46%  '
$c = $ws.Range('D22')
$c.Value = '''371.77'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -5.11%  '
$ws.Range('D23').Value = '3.75'
$c = $ws.Range('D24')
$c.Value = '''80.22'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -3.40%  '
$ws.Range('D25').Value = '10.99'
$ws.Range('E25').Value = '  -9.77%  '
$c = $ws.Range('D26')
$c.Value = '''3.83'
$c.Style = "Normal"
$ws.Range('E26').Value = '  -2.56%  '
$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D27')
$c.Value = '''2.65'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -5.75%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '11.31'
$ws.Range('E28').Value = '  -7.18%  '
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D29')
$c.Value = '''8.28'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -6.19%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D30')
$c.Value = '''28.65'
$c.Style = "Normal"
$ws.Range('E30').Value = '  -7.44%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range('D31')
$c.Value = '''634.28'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -8.40%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D32')
$c.Value = '''6.58'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -6.20%  '
$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D33').Value = '11.17'
$ws.Range('E33').Value = '  -4.52%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.105'
$ws.Range('E34').Value = '  -5.88%  '
$ws.Range('D35').Value = '58.75'
$ws.Range('E35').Value = '  -7.87%  '
$ws.Range('B36').Value = 'Dai'
$ws.Range('C36').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D36')
$c.Value = '''1.00'
$c.Style = "Normal"
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('B37').Value = 'TheGraph'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Range('D37')
$c.Value = '''0.386'
$c.Style = "Normal"
$ws.Range('E37').Value = '  -2.31%  '
$ws.Range('B38').Value = 'InjectiveProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range('D38')
$c.Value = '''36.33'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -5.86%  '
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  +0.10%  '
$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D40').Value = '0.0₃0699'
$ws.Range('E40').Value = '  +3.18%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '0.123'
$ws.Range('E41').Value = '  -5.31%  '
$ws.Range('D42').Value = '2.910.27'
$ws.Range('E42').Value = '  -4.83%  '
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range('D43')
$c.Value = '''2.44'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -4.99%  '
$ws.Range('B44').Value = 'ThetaToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$c = $ws.Range('D44')
$c.Value = '''2.67'
$c.Style = "Normal"
$ws.Range('E44').Value = '  -10.42%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range('D45')
$c.Value = '''2.64'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -4.08%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D46')
$c.Value = '''0.0394'
$c.Style = "Normal"
$ws.Range('E46').Value = '  -2.30%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c = $ws.Range('D47')
$c.Value = '''2.99'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +3.64%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range('D48')
$c.Value = '''2.76'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +5.44%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = '0.125'
$ws.Range('E49').Value = '  -2.00%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D50')
$c.Value = '''134.96'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -2.45%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range('D51')
$c.Value = '''2.39'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -10.23%  '
